# Implemented the CHAR function in C and Ruby
# This script adds a new "LEN(B)" column (C) for rows 15-47, and adds
# several new edge-case rows (271-276) exercising CHAR()/error handling
# at the bottom of the "character table" on sheet 2 (工作表1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# New column C: LEN(B#) for every row that currently has a CHAR() result
# in column B (rows 15 through 47). Rows 15-45 form one contiguous shared
# formula block; 46 and 47 follow the same pattern as the source file.
$ws.Range("C15:C45").Formula = "=LEN(B15)"
$ws.Range("C46").Formula = "=LEN(B46)"
$ws.Range("C47").Formula = "=LEN(B47)"

# Extend the existing B208:B271 shared CHAR() formula down to the new
# row 271 by re-applying it across the full (now-longer) column range.
$ws.Range("B208:B271").Formula = "=CHAR(A208)"

# Row 271: 0/0 division error feeding into CHAR()
$ws.Range("A271").Formula = "=0/0"

# Row 272: CHAR() of a blank cell (A272 left empty) -> #VALUE!
$ws.Range("B272").Formula = "=CHAR(A272)"

# Row 273: A273 is the text "48" (not numeric 48) -> CHAR() on a numeric
# string coerces and evaluates normally.
$ws.Range("A273").Formula = "=""48"""
$ws.Range("B273:B276").Formula = "=CHAR(A273)"

# Row 274: A274 is an ordinary text string -> CHAR() on text -> #VALUE!
$ws.Range("A274").Value = "asfasd"

# Row 275: A275 is a negative number -> CHAR() on negative -> #VALUE!
$ws.Range("A275").Value = -1

# Row 276: A276 is a non-integer number that truncates to a valid code
# point (97 -> "a").
$ws.Range("A276").Value = 97.12312

# Update the visible window: scroll near the top of the table and select
# the newly added LEN() column.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C15:C45").Select()
